# Actualizacion plan de calidad tiempos de ejecucion
#
# Changes applied:
#  1. C27 and C28 ("Garantia" section, "Momento de ejecucion" column) change
#     from "Mensual" to "Al finalizar proyecto".
#  2. A new (empty) cell C29 is created with an underline font style.
#  3. The active selection moves to C29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the execution-moment values for the warranty rows.
$ws.Range("C27").Value = "Al finalizar proyecto"
$ws.Range("C28").Value = "Al finalizar proyecto"

# 2. Create the new styled (underlined) empty cell at C29.
$c29 = $ws.Range("C29")
$c29.Value = ""
$c29.Font.Underline = 2

# 3. Move the active selection to C29, matching the saved view state.
$c29.Select() | Out-Null
